$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates to whole numbers.
$ws.Range("Q2").Value = 776795
$ws.Range("R2").Value = 7198204

# Drop the (empty/placeholder) start- and end-time values; the date
# columns (Y2/AA2) are left untouched.
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
